$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G3").Value = "2016-08-15 20:42:43"

$wsZhCn.Range("H3").Value = "2016-08-15 20:42:39"
$wsZhCn.Range("K3").Value = "2016-08-15 20:42:56"

$wsDeDe.Range("H3").Value = "2016-08-15 20:42:43"
$wsDeDe.Range("K3").Value = "2016-08-15 20:43:13"
